# Generate Report for Handback
# Updates the localization-status workbook to reflect a completed handback:
#  - Overview sheet: zh-cn / de-de status columns move from "Ready for
#    handoff" to "Handed back: in sync with en-US"
#  - zh-cn / de-de detail sheets: Latest Target File / Latest Handback File
#    columns get populated (with hyperlinks on the target file) and the
#    de-de sheet's Latest Handback DateTime is stamped with the handback
#    time.
#  - A few columns are widened so the new, longer strings aren't truncated.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"
$hyperlinkColor = 15570276   # BGR-encoded 0x6495ED (CornflowerBlue) -> matches existing custom "HyperLink" style

$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/bfc158a46f44c05954ddf1e4351a46a8b60b93f3/e2e/"
$mdName1 = "346cf65b-951e-4cf7-b8c1-8e7fe9775cb8.md"
$mdName2 = "fe03d05e-fd73-4c39-a5c4-bac138d75e85.md"

# ---------------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")

$overview.Range("E2").Value = $newStatus
$overview.Range("F2").Value = $newStatus
$overview.Range("E3").Value = $newStatus
$overview.Range("F3").Value = $newStatus

$overview.Columns.Item(5).ColumnWidth = 29.1667
$overview.Columns.Item(6).ColumnWidth = 29.1667

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("C2").Value = $newStatus
$zhcn.Range("C3").Value = $newStatus

# Latest Target File (I) / Latest Handback File (J) for each source file
$zhcn.Range("I2").Value = $mdName1
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), ($baseUrl + $mdName1), "", "", $mdName1)
$zhcn.Range("I2").Font.Underline = 2
$zhcn.Range("I2").Font.Color = $hyperlinkColor
$zhcn.Range("J2").Value = "346cf65b-951e-4cf7-b8c1-8e7fe9775cb8.b5b2c78f74a18d86b69be7ee640f3f5a6e766176.zh-cn.xlf"

$zhcn.Range("I3").Value = $mdName2
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), ($baseUrl + $mdName2), "", "", $mdName2)
$zhcn.Range("I3").Font.Underline = 2
$zhcn.Range("I3").Font.Color = $hyperlinkColor
$zhcn.Range("J3").Value = "fe03d05e-fd73-4c39-a5c4-bac138d75e85.2013f16f08f420be2005c28307c7114c24aff14f.zh-cn.xlf"

# Latest Handback DateTime was never set ("0001-01-01 00:00:00"); this
# handback run stamps it with the same generation time already recorded
# elsewhere for this run.
$zhcn.Range("K2").Value = "2016-08-17 04:42:57"
$zhcn.Range("K3").Value = "2016-08-17 04:42:57"

$zhcn.Columns.Item(3).ColumnWidth = 29.1667
$zhcn.Columns.Item(9).ColumnWidth = 39.1667
$zhcn.Columns.Item(10).ColumnWidth = 39.1667

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("C2").Value = $newStatus
$dede.Range("C3").Value = $newStatus

$dede.Range("I2").Value = $mdName1
$dede.Hyperlinks.Add($dede.Range("I2"), ($baseUrl + $mdName1), "", "", $mdName1)
$dede.Range("I2").Font.Underline = 2
$dede.Range("I2").Font.Color = $hyperlinkColor
$dede.Range("J2").Value = "346cf65b-951e-4cf7-b8c1-8e7fe9775cb8.b5b2c78f74a18d86b69be7ee640f3f5a6e766176.de-de.xlf"

$dede.Range("I3").Value = $mdName2
$dede.Hyperlinks.Add($dede.Range("I3"), ($baseUrl + $mdName2), "", "", $mdName2)
$dede.Range("I3").Font.Underline = 2
$dede.Range("I3").Font.Color = $hyperlinkColor
$dede.Range("J3").Value = "fe03d05e-fd73-4c39-a5c4-bac138d75e85.2013f16f08f420be2005c28307c7114c24aff14f.de-de.xlf"

# de-de's handback finished a little later than zh-cn's, so it gets its
# own (later) handback timestamp.
$dede.Range("K2").Value = "2016-08-17 04:43:09"
$dede.Range("K3").Value = "2016-08-17 04:43:09"

$dede.Columns.Item(3).ColumnWidth = 29.1667
$dede.Columns.Item(9).ColumnWidth = 39.1667
$dede.Columns.Item(10).ColumnWidth = 39.1667
